$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")

# Clear the old value in Hoja2!C6 ("sdfas") and move content to C1
$ws2.Range("C6").Clear()

# New formula cell referencing Hoja1!A1
$ws2.Range("C1").Formula = "=Hoja1!A1"

# Highlight the new formula cell with a yellow fill
$ws2.Range("C1").Interior.Color = 65535

# Update selection / active cell on Hoja2 to C1 without changing the
# workbook's active sheet (Hoja1 stays the selected tab, as before).
$ws2.Activate()
$ws2.Range("C1").Select()
$ws1.Activate()

$wb.Application.CalculateFull()
